$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New API row: "Update Chat Bot Agent Status (PUT)"
$url = "http://1msg.1point1.in:3001/api/chat/bot/update/agent/status/"

$ws.Range("A48").Value = 19
$ws.Range("B48").Value = "Update Chat Bot  Agent Status  ( PUT )"
$ws.Range("C48").Value = "PUT"
$ws.Range("D48").Value = $url
$ws.Range("F48").Value = "{`n    ""user_id"" : 10,`n    ""agent_id"" : 20,`n    ""agent_status"" : true`n}"
$ws.Range("G48").Value = "curl --location --request PUT 'http://1msg.1point1.in:3001/api/chat/bot/update/agent/status/' \`n--header 'Content-Type: application/json' \`n--data '{`n    ""user_id"" : 10,`n    ""agent_id"" : 20,`n    ""agent_status"" : true`n}'`n"

# Match the row height used by the other full-detail rows (e.g. row 45)
$ws.Rows.Item(48).RowHeight = 187.2

# Hyperlink the URL cell, then restore the wrapped hyperlink formatting used
# elsewhere in the sheet (Hyperlinks.Add overwrites the cell style, so copy
# the formatting from an existing hyperlinked cell afterwards).
$ws.Hyperlinks.Add($ws.Range("D48"), $url) | Out-Null
$ws.Range("D45").Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update the view state to match: scrolled down with G48 selected.
$ws.Range("G48").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 2
